$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look numeric must be forced to Text format
# so Excel stores them as literal strings (matching the source feeds formatting),
# not auto-converted numbers.
$textCells = @("D4", "D5", "D6", "D7", "D8", "D15", "D16", "D18", "D20", "D21", "D25", "D27", "D30", "D32", "D36", "D40", "D41", "D43", "D44", "D46", "D47", "D48", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.975.57"
$ws.Range("E2").Value = "  +3.38%  "
$ws.Range("D3").Value = "1.726.95"
$ws.Range("E3").Value = "  +3.10%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "218.58"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "0.525"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "24.07"
$ws.Range("E8").Value = "  +13.59%  "
$ws.Range("E9").Value = "  +3.57%  "
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").Value = "1.970.89"
$ws.Range("E12").Value = "  +3.12%  "
$ws.Range("D13").Value = "1.724.26"
$ws.Range("E13").Value = "  +2.86%  "
$ws.Range("E14").Value = "  +3.66%  "
$ws.Range("D15").Value = "0.567"
$ws.Range("E15").Value = "  +6.26%  "
$ws.Range("D16").Value = "67.91"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "27.923.67"
$ws.Range("D18").Value = "243.38"
$ws.Range("E18").Value = "  +2.71%  "
$ws.Range("D19").Value = "0.0₃0759"
$ws.Range("E19").Value = "  +2.52%  "
$ws.Range("D20").Value = "7.89"
$ws.Range("E20").Value = "  -3.11%  "
$ws.Range("D21").Value = "1.00"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +4.56%  "
$ws.Range("E23").Value = "  +4.77%  "
$ws.Range("E24").Value = "  +0.50%  "
$ws.Range("D25").Value = "149.62"
$ws.Range("E26").Value = "  +4.37%  "
$ws.Range("D27").Value = "16.83"
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("E28").Value = "  +1.87%  "
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").Value = "0.0511"
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").Value = "3.46"
$ws.Range("E32").Value = "  +2.94%  "
$ws.Range("E33").Value = "  +3.34%  "
$ws.Range("D34").Value = "1.483.56"
$ws.Range("E34").Value = "  -3.58%  "
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("D36").Value = "0.612"
$ws.Range("E36").Value = "  +2.68%  "
$ws.Range("E37").Value = "  +4.81%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("D40").Value = "1.08"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Value = "71.54"
$ws.Range("E41").Value = "  +5.50%  "
$ws.Range("E42").Value = "  +5.03%  "
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.15%  "
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  +1.66%  "
$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.874.75"
$ws.Range("E45").Value = "  +3.14%  "
$ws.Range("D46").Value = "0.791"
$ws.Range("E46").Value = "  +0.85%  "
$ws.Range("D47").Value = "1.77"
$ws.Range("E47").Value = "  +13.44%  "
$ws.Range("D48").Value = "91.55"
$ws.Range("E48").Value = "  +0.76%  "
$ws.Range("E49").Value = "  +3.80%  "
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").Value = "8.23"
$ws.Range("E51").Value = "  +2.29%  "
